$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 181, shifting existing rows 181:256 down to 182:256
$ws.Rows.Item(181).Insert()

# Populate the new row 181 with fresh data
$ws.Cells.Item(181, 1).Value = 10
$ws.Cells.Item(181, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(181, 3).Value = "La Araucanía"
$ws.Cells.Item(181, 4).Value = 44523
$ws.Cells.Item(181, 5).Value = 9
$ws.Cells.Item(181, 6).Value = 100112037
$ws.Cells.Item(181, 7).Value = "Cebollín"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 30
$ws.Cells.Item(181, 11).Value = 9000
$ws.Cells.Item(181, 12).Value = 9000
$ws.Cells.Item(181, 13).Value = 9000
$ws.Cells.Item(181, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(181, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(181, 16).Value = 750
$ws.Cells.Item(181, 17).Value = 12
$ws.Cells.Item(181, 18).Value = "Hortaliza"
